# City Bike project.docx - proofreading fix-up
# Removes the typo/misspellings (and their now-stale w:proofErr squiggle
# markers) that were present in the original draft, per the commit diff.

$d = $word.ActiveDocument

function Fix($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# "Our task ..." paragraph -----------------------------------------------
Fix "Logs to builda dashboard" "Logs to build a dashboard"
Fix "work with a timespam of" "work with a timespan of"
Fix "periods. The folloowing are" "periods. The following are"
Fix "are the quetions you may" "are the questions you may"

# "Task need to to do achieved:" -----------------------------------------
Fix "Task need to to do achieved:" "Task need to do achieved:"

# "Acquire data ..." -------------------------------------------------------
Fix "some of the datasers(SQL, MongoDB, casendra local" "some of the datasets (SQL, MongoDB, Casandra local"

# "Connect with the business users and try to get the understanding about KPI(..." 
Fix "understanding about KPI(Key performance Indicator)" "understanding about KPI (Key performance Indicator)"

# "A KPI is a measureable value ..." --------------------------------------
Fix "A KPI is a measureable value" "A KPI is a measurable value"
Fix "demonstrates how effecivly a company" "demonstrates how effectively a company"
Fix "achieving key business abjective." "achieving key business objective."

# "Connect with business user with raw visualization ..." -----------------
Fix "expectations feedback baced on ese of use." "expectations feedback based on ese of use."

# "Decise total number of dashboard ..." ----------------------------------
Fix "Decise total number of dashboard based" "Decide total number of dashboards based"

# "Start building production-based dashboaed." ----------------------------
Fix "Start building production-based dashboaed." "Start building production-based dashboard."

# Pick hours during summer/winter -----------------------------------------
Fix "bikes are used dusring summer month?" "bikes are used during summer month?"
Fix "bikes are used dusring winter month?" "bikes are used during winter month?"

# Top/bottom 10 stations questions ----------------------------------------
Fix "starting a journey?(Based on data, why do you hypothesize these are the top locations?)" `
    "starting a journey? (Based on data, why do you hypothesize these are the top locations?)"
Fix "ending a journey?(Based on data, why?)" "ending a journey? (Based on data, why?)"
Fix "starting a journey?(Based on data, why?)" "starting a journey? (Based on data, why?)"
Fix "ending a journey?(Based on data, why ?)" "ending a journey? (Based on data, why?)"

# Gender breakdown question -------------------------------------------------
Fix "active partcipants(Male vs Female)?" "active participants (Male vs Female)?"

# Gender outreach question --------------------------------------------------
Fix "How effective has gender outrach been" "How effective has gender outreach been"
